# Updated cryptos list values (price + 1h volume change) per target diff.
# Cells are stored as text in the workbook (inline strings), so we force
# NumberFormat to "@" (Text) before assigning, which prevents Excel's COM
# layer from auto-converting numeric-looking strings (e.g. "1.002",
# "0.000006754", "3.460") into actual numbers / scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.063.16'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.763.88'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.06'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5232'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2758'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.43'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06204'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.777.75'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07014'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6419'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.543'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '78.25'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '26.088.24'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.67'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006754'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.001.58'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.079'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.453'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.204'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.08'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.514'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.98%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '103.25'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08405'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.707'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.460'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04468'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.619'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6072'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.745'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01594'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.990'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.33%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.80'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3888'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7450'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.945'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05514'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.371'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.37%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.27'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.76'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.33%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.83%  '
